# Apply updated TPM values + relabeled MuSCs/Inflammatory-Mac rows, and 4 new rows (10-13)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna3"
$ws.Range("C2").Value = "Epha4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2420556666666667
$ws.Range("H2").Value = 0.726167
$ws.Range("I2").Value = 0.5314769098578004
$ws.Range("J2").Value = 0.5314769098578004
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.540560666666666
$ws.Range("N2").Value = 25.621682
$ws.Range("O2").Value = 0.4159358086620884
$ws.Range("P2").Value = 0.4159358086620884
$ws.Range("Q2").Value = 2.067291105877111
$ws.Range("R2").Value = 18.605619952894
$ws.Range("S2").Value = 0.2210602782869321
$ws.Range("T2").Value = 0.2210602782869321

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna3"
$ws.Range("C3").Value = "Epha4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2420556666666667
$ws.Range("H3").Value = 0.726167
$ws.Range("I3").Value = 0.5314769098578004
$ws.Range("J3").Value = 0.5314769098578004
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.57455166666667
$ws.Range("N3").Value = 34.723655
$ws.Range("O3").Value = 0.563694901924408
$ws.Range("P3").Value = 0.563694901924408
$ws.Range("Q3").Value = 2.801685820042778
$ws.Range("R3").Value = 25.215172380385
$ws.Range("S3").Value = 0.2995908245773802
$ws.Range("T3").Value = 0.2995908245773802

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna3"
$ws.Range("C4").Value = "Epha4"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2420556666666667
$ws.Range("H4").Value = 0.726167
$ws.Range("I4").Value = 0.5314769098578004
$ws.Range("J4").Value = 0.5314769098578004
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.026642
$ws.Range("N4").Value = 0.079926
$ws.Range("O4").Value = 0.001297498167494471
$ws.Range("P4").Value = 0.001297498167494471
$ws.Range("Q4").Value = 0.006448847071333333
$ws.Range("R4").Value = 0.058039623642
$ws.Range("S4").Value = 0.0006895903166061202
$ws.Range("T4").Value = 0.0006895903166061202

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna3"
$ws.Range("C5").Value = "Epha4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2420556666666667
$ws.Range("H5").Value = 0.726167
$ws.Range("I5").Value = 0.5314769098578004
$ws.Range("J5").Value = 0.5314769098578004
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.391608
$ws.Range("N5").Value = 1.174824
$ws.Range("O5").Value = 0.01907179124600912
$ws.Range("P5").Value = 0.01907179124600912
$ws.Range("Q5").Value = 0.094790935512
$ws.Range("R5").Value = 0.8531184196080001
$ws.Range("S5").Value = 0.01013621667688197
$ws.Range("T5").Value = 0.01013621667688197

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna3"
$ws.Range("C6").Value = "Epha4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2054156666666667
$ws.Range("H6").Value = 0.616247
$ws.Range("I6").Value = 0.4510271759376837
$ws.Range("J6").Value = 0.4510271759376837
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 8.540560666666666
$ws.Range("N6").Value = 25.621682
$ws.Range("O6").Value = 0.4159358086620884
$ws.Range("P6").Value = 0.4159358086620884
$ws.Range("Q6").Value = 1.754364963050444
$ws.Range("R6").Value = 15.789284667454
$ws.Range("S6").Value = 0.1875983531522185
$ws.Range("T6").Value = 0.1875983531522185

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna3"
$ws.Range("C7").Value = "Epha4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2054156666666667
$ws.Range("H7").Value = 0.616247
$ws.Range("I7").Value = 0.4510271759376837
$ws.Range("J7").Value = 0.4510271759376837
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 11.57455166666667
$ws.Range("N7").Value = 34.723655
$ws.Range("O7").Value = 0.563694901924408
$ws.Range("P7").Value = 0.563694901924408
$ws.Range("Q7").Value = 2.377594246976111
$ws.Range("R7").Value = 21.398348222785
$ws.Range("S7").Value = 0.2542417197054353
$ws.Range("T7").Value = 0.2542417197054353

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna3"
$ws.Range("C8").Value = "Epha4"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2054156666666667
$ws.Range("H8").Value = 0.616247
$ws.Range("I8").Value = 0.4510271759376837
$ws.Range("J8").Value = 0.4510271759376837
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.026642
$ws.Range("N8").Value = 0.079926
$ws.Range("O8").Value = 0.001297498167494471
$ws.Range("P8").Value = 0.001297498167494471
$ws.Range("Q8").Value = 0.005472684191333333
$ws.Range("R8").Value = 0.04925415772199999
$ws.Range("S8").Value = 0.000585206934269351
$ws.Range("T8").Value = 0.000585206934269351

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna3"
$ws.Range("C9").Value = "Epha4"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2054156666666667
$ws.Range("H9").Value = 0.616247
$ws.Range("I9").Value = 0.4510271759376837
$ws.Range("J9").Value = 0.4510271759376837
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.391608
$ws.Range("N9").Value = 1.174824
$ws.Range("O9").Value = 0.01907179124600912
$ws.Range("P9").Value = 0.01907179124600912
$ws.Range("Q9").Value = 0.080442418392
$ws.Range("R9").Value = 0.7239817655280001
$ws.Range("S9").Value = 0.008601896145760528
$ws.Range("T9").Value = 0.008601896145760528

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Efna3"
$ws.Range("C10").Value = "Epha4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.007968333333333332
$ws.Range("H10").Value = 0.023905
$ws.Range("I10").Value = 0.01749591420451593
$ws.Range("J10").Value = 0.01749591420451593
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 8.540560666666666
$ws.Range("N10").Value = 25.621682
$ws.Range("O10").Value = 0.4159358086620884
$ws.Range("P10").Value = 0.4159358086620884
$ws.Range("Q10").Value = 0.06805403424555555
$ws.Range("R10").Value = 0.6124863082099999
$ws.Range("S10").Value = 0.007277177222937851
$ws.Range("T10").Value = 0.007277177222937851

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Efna3"
$ws.Range("C11").Value = "Epha4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.007968333333333332
$ws.Range("H11").Value = 0.023905
$ws.Range("I11").Value = 0.01749591420451593
$ws.Range("J11").Value = 0.01749591420451593
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 11.57455166666667
$ws.Range("N11").Value = 34.723655
$ws.Range("O11").Value = 0.563694901924408
$ws.Range("P11").Value = 0.563694901924408
$ws.Range("Q11").Value = 0.09222988586388887
$ws.Range("R11").Value = 0.830068972775
$ws.Range("S11").Value = 0.009862357641592462
$ws.Range("T11").Value = 0.009862357641592462

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Efna3"
$ws.Range("C12").Value = "Epha4"
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.007968333333333332
$ws.Range("H12").Value = 0.023905
$ws.Range("I12").Value = 0.01749591420451593
$ws.Range("J12").Value = 0.01749591420451593
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.026642
$ws.Range("N12").Value = 0.079926
$ws.Range("O12").Value = 0.001297498167494471
$ws.Range("P12").Value = 0.001297498167494471
$ws.Range("Q12").Value = 0.0002122923366666666
$ws.Range("R12").Value = 0.00191063103
$ws.Range("S12").Value = 0.0000227009166189999
$ws.Range("T12").Value = 0.0000227009166189999

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Efna3"
$ws.Range("C13").Value = "Epha4"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.007968333333333332
$ws.Range("H13").Value = 0.023905
$ws.Range("I13").Value = 0.01749591420451593
$ws.Range("J13").Value = 0.01749591420451593
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.391608
$ws.Range("N13").Value = 1.174824
$ws.Range("O13").Value = 0.01907179124600912
$ws.Range("P13").Value = 0.01907179124600912
$ws.Range("Q13").Value = 0.00312046308
$ws.Range("R13").Value = 0.02808416772
$ws.Range("S13").Value = 0.0003336784233666134
$ws.Range("T13").Value = 0.0003336784233666134

